$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style of an existing header cell (E1) to the new headers so they match (bold, centered, bordered)
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

$boolValues = @(
    @(0,0,0),
    @(0,0,0),
    @(1,0,0),
    @(1,1,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0)
)

for ($i = 0; $i -lt $boolValues.Length; $i++) {
    $row = $i + 2
    $vals = $boolValues[$i]
    $ws.Cells.Item($row, 6).Value = [bool]($vals[0] -eq 1)
    $ws.Cells.Item($row, 7).Value = [bool]($vals[1] -eq 1)
    $ws.Cells.Item($row, 8).Value = [bool]($vals[2] -eq 1)
}
